$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Delete()

$ws.Range("A1").Value2 = "test@gmail.com"
$ws.Range("B1").Value2 = "test"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:nbanish@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:Thodupuzha@1")

Write-Output "done"
